# Insert new data in the rows of the table on the "a" (SubjectWithGradesView) sheet.
# Columns: A=year, B=semester, C=course code, D=course title, E=units, F=grades

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("a")

# New rows to append right after the existing last row (78).
$newRows = @(
    @("1", "1", "test1", "testing", "3.0", "not yet taken"),
    @("1", "2", "test2", "testing", "2.0", "69"),
    @("3", "2", "test3", "test",    "3.0", "79"),
    @("4", "2", "test4", "tesss",   "2.0", "99"),
    @("3", "2", "test5", "tess",    "2.0", "55")
)

$columns = @("A", "B", "C", "D", "E", "F")
$startRow = 79

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $newRows[$i]

    for ($j = 0; $j -lt $columns.Count; $j++) {
        $cellRef = "$($columns[$j])$rowNum"
        $value = $rowValues[$j]
        $cell = $ws.Range($cellRef)

        if ($value -match '^[0-9]+(\.[0-9]+)?$') {
            # Value looks numeric (e.g. "1", "2.0", "69") but must be stored as
            # literal text, matching the existing rows above it. Enter it with
            # a leading apostrophe (forces text) then strip the resulting
            # quote-prefix formatting so no new cell style is introduced.
            $cell.Formula = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}

Write-Host "Inserted rows 79-83 on sheet 'a'"
